$d = $word.ActiveDocument

# Pull the whole package as flat OOXML so we can perform part-level
# surgery (removing the header part + its relationship) that isn't
# reachable through the regular object model (HeaderFooter.Exists is
# read-only in Word's COM model, same as in real Word).
$xml = $d.Content.WordOpenXML

# 1) Drop the two leading empty paragraphs at the start of the body.
$xml = $xml.Replace('<w:body><w:p w14:paraId="510ED0D8" w14:textId="2D9707AC" w:rsidR="001230EA" w:rsidRDefault="001230EA"/><w:p w14:paraId="202C1F4E" w14:textId="77777777" w:rsidR="00936D2F" w:rsidRDefault="00936D2F"/><w:tbl>', '<w:body><w:tbl>')

# 2) Remove the header reference from the (only) section properties.
$xml = $xml.Replace('<w:headerReference w:type="default" r:id="rId6"/>', '')

# 3) Remove the header1.xml part entirely (package part + its XML data).
$xml = $xml -replace '<pkg:part pkg:name="/word/header1\.xml"[^>]*><pkg:xmlData>.*?</pkg:xmlData></pkg:part>', ''

# 4) Remove the now-dangling relationship pointing at header1.xml.
$xml = $xml.Replace('<Relationship Id="rId6" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/header" Target="header1.xml"/>', '')

$d.Content.InsertXML($xml)
